# add save column in s_vals sheets
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the "sum" header's formatting (bold font + border, style index 1)
# onto the new "Save" header cell, then set its text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New data column, defaulting to 0 for each existing row.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
